$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain text, matching original inlineStr formatting
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.400.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.603.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.617.31'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.076.38'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.213.10'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.547.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.169'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.405'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.718.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0742'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.79'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '149.77'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.865'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.868'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.08'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '270.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.68'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0955'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.966.97'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.60'
$ws.Range('D51').Style = 'Normal'

# Other text-only cells (names, links, percentages)
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('E16').Value = '  -2.39%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('E22').Value = '  +3.30%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E30').Value = '  -5.51%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  -6.20%  '
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('E37').Value = '  -3.83%  '
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('E40').Value = '  +2.67%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('E51').Value = '  +2.16%  '
